$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.527.25'
$ws.Range("E2").Value = '  +5.57%  '

$ws.Range("D3").Value = '2.054.56'
$ws.Range("E3").Value = '  +3.99%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").Value = '252.41'
$ws.Range("E5").Value = '  +2.93%  '

$ws.Range("D6").Value = '0.650'
$ws.Range("E6").Value = '  +3.09%  '

$ws.Range("D7").Value = '66.32'
$ws.Range("E7").Value = '  +15.65%  '

$ws.Range("E8").Value = '  +0.02%  '

$ws.Range("E9").Value = '  +6.44%  '

$ws.Range("D10").Value = '59.97'
$ws.Range("E10").Value = '  +2.95%  '

$ws.Range("E11").Value = '  +4.71%  '

$ws.Range("E12").Value = '  +1.30%  '

$ws.Range("D13").Value = '0.907'
$ws.Range("E13").Value = '  -3.95%  '

$ws.Range("D14").Value = '14.95'
$ws.Range("E14").Value = '  +4.93%  '

$ws.Range("D15").Value = '2.357.29'
$ws.Range("E15").Value = '  +4.19%  '

$ws.Range("D16").Value = '21.38'
$ws.Range("E16").Value = '  +21.74%  '

$ws.Range("D17").Value = '5.58'
$ws.Range("E17").Value = '  +6.18%  '

$ws.Range("D18").Value = '2.054.94'
$ws.Range("E18").Value = '  +3.82%  '

$ws.Range("D19").Value = '37.331.00'
$ws.Range("E19").Value = '  +5.43%  '

$ws.Range("D20").Value = '73.73'
$ws.Range("E20").Value = '  +3.20%  '

$ws.Range("D21").Value = '0.0₃0877'
$ws.Range("E21").Value = '  +4.43%  '

$ws.Range("D22").Value = '5.45'
$ws.Range("E22").Value = '  +6.27%  '

$ws.Range("D23").Value = '240.25'
$ws.Range("E23").Value = '  +3.35%  '

$ws.Range("E24").Value = '  +2.26%  '

$ws.Range("E25").Value = '  -0.08%  '

$ws.Range("D26").Value = '2.39'
$ws.Range("E26").Value = '  +4.07%  '

$ws.Range("E27").Value = '  +8.29%  '

$ws.Range("D28").Value = '160.39'
$ws.Range("E28").Value = '  -1.98%  '

$ws.Range("E29").Value = '  +4.39%  '

$ws.Range("D30").Value = '5.26'
$ws.Range("E30").Value = '  +8.19%  '

$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D31").Value = '0.122'
$ws.Range("E31").Value = '  +3.25%  '

$ws.Range("B32").Value = 'Kaspa'
$ws.Range("C32").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D32").Value = '0.114'
$ws.Range("E32").Value = '  +24.50%  '

$ws.Range("E33").Value = '  +6.94%  '

$ws.Range("D34").Value = '4.76'
$ws.Range("E34").Value = '  +11.53%  '

$ws.Range("E35").Value = '  +5.40%  '

$ws.Range("E36").Value = '  +3.86%  '

$ws.Range("B37").Value = 'WEMIXToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D37").Value = '1.85'
$ws.Range("E37").Value = '  +4.23%  '

$ws.Range("B38").Value = 'BinanceUSD'
$ws.Range("C38").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D38").Value = '1.00'
$ws.Range("E38").Value = '  -0.05%  '

$ws.Range("D39").Value = '6.13'
$ws.Range("E39").Value = '  +20.04%  '

$ws.Range("E40").Value = '  +34.43%  '

$ws.Range("E41").Value = '  +16.61%  '

$ws.Range("E42").Value = '  +3.04%  '

$ws.Range("E43").Value = '  +4.82%  '

$ws.Range("D44").Value = '0.0219'
$ws.Range("E44").Value = '  +3.94%  '

$ws.Range("E45").Value = '  +6.16%  '

$ws.Range("D46").Value = '16.96'
$ws.Range("E46").Value = '  +6.90%  '

$ws.Range("B47").Value = 'FraxShare'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D47").Value = '7.98'
$ws.Range("E47").Value = '  +6.63%  '

$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").Value = '95.58'
$ws.Range("E48").Value = '  +4.77%  '

$ws.Range("D49").Value = '1.420.99'
$ws.Range("E49").Value = '  +3.04%  '

$ws.Range("E50").Value = '  +2.37%  '

$ws.Range("D51").Value = '46.57'
$ws.Range("E51").Value = '  +1.68%  '
